# Add three new VIC Mystery Cases rows at the top of the data table:
#   2020-12-30 | 3132 | (new DHHS link)
#   2020-12-30 | 3803 | (new DHHS link)
#   2020-12-27 | 3194 | (new DHHS link)
# The rest of the existing rows shift down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right below the header (row 1), pushing the
# existing data down.
$ws.Range("A2:E4").EntireRow.Insert() | Out-Null

# The inserted rows pick up the header's formatting by default; copy the
# number/font formatting from the row just below (the original first
# data row, now shifted to row 5) onto the 3 new rows.
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A2:E4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Grow the table / AutoFilter so it covers the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E28"))

$newLink = "https://www.dhhs.vic.gov.au/victoria-records-three-new-community-coronavirus-cases-30-december-2020"

$newRows = @(
    @{ Row = 2; Date = "2020-12-30"; Postcode = 3132 },
    @{ Row = 3; Date = "2020-12-30"; Postcode = 3803 },
    @{ Row = 4; Date = "2020-12-27"; Postcode = 3194 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.Date
    $ws.Cells.Item($r, 2).Value = $row.Postcode
    $ws.Cells.Item($r, 3).Value = $newLink
    $ws.Cells.Item($r, 4).Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
    $ws.Cells.Item($r, 5).Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"
}

# Resizing the table rewrote the calculated-column formulas of the rows
# that were newly absorbed at the bottom (26-28) into a form that doesn't
# evaluate correctly - restore them explicitly.
foreach ($r in 26..28) {
    $ws.Cells.Item($r, 4).Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
    $ws.Cells.Item($r, 5).Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"
}

# The hyperlink in column C moved with the data (old row 7 -> new row
# 10); the stored hyperlink reference doesn't shift automatically, so
# re-home it manually.
$ws.Range("C7").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-october-2020") | Out-Null
# Adding the hyperlink re-styles the cell; restore the plain data style
# used by the rest of the column.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B4").Select()
